$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Files" query (row 4, column B) is corrected: the `File Type` column
# and the `Breed` column are removed from the returned projection.
$ws.Range("B4").Value = 'MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN[''Greyhound''] 
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '''') AS `File Name`,
         coalesce(labels(parent)[0], '''') AS `Association`,
        coalesce(f.file_description, '''') AS `Description`,
        coalesce(f.file_format, '''') AS `Format`,
        coalesce(f.file_size, '''') AS `Size`,
        coalesce(c.case_id, '''') AS `Case ID`,
         coalesce(diag.disease_term,'''') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'''') AS `Study Code`'

# Row 4 shrank (two fewer lines in the wrapped query text), so its
# auto-fitted height drops from 246.5 to 217.5.
$ws.Rows.Item(4).RowHeight = 217.5

# The selection / active cell moved from B2 down to B4.
$ws.Range("B4").Select()
